# Auto-generated edit script: append rows 214-219 (dates 46007-46009,
# stations 四方坪站/高岭站) to Sheet1, then reselect D226 to match the
# commit's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 214: date serial 46007, station string-id 26
$ws.Cells.Item(214, 1).Value = 46007
$ws.Cells.Item(214, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(214, 3).Value = 626.23399999999992
$ws.Cells.Item(214, 4).Value = 1225.2180000000003
$ws.Cells.Item(214, 5).Value = 394.75000000000006
$ws.Cells.Item(214, 6).Value = 327.44
$ws.Cells.Item(214, 7).Value = 96.724999999999994
$ws.Cells.Item(214, 8).Value = 536.00300000000004
$ws.Cells.Item(214, 9).Value = 304.06199999999995
$ws.Cells.Item(214, 10).Value = 90.1
$ws.Cells.Item(214, 11).Value = 96.079999999999984
$ws.Cells.Item(214, 12).Value = 54.74
$ws.Cells.Item(214, 13).Value = 239.29599999999996
$ws.Cells.Item(214, 14).Value = 299.88800000000003
$ws.Cells.Item(214, 15).Value = 583.71799999999996
$ws.Cells.Item(214, 16).Value = 1265.8440000000001
$ws.Cells.Item(214, 17).Value = 476.36300000000006
$ws.Cells.Item(214, 18).Value = 420.721
$ws.Cells.Item(214, 19).Value = 190.25599999999997
$ws.Cells.Item(214, 20).Value = 181.13799999999998
$ws.Cells.Item(214, 21).Value = 83.837000000000003
$ws.Cells.Item(214, 22).Value = 145.26
$ws.Cells.Item(214, 23).Value = 138.06
$ws.Cells.Item(214, 24).Value = 34.488999999999997
$ws.Cells.Item(214, 25).Value = 59.3
$ws.Cells.Item(214, 26).Value = 105.56

# Row 215: date serial 46007, station string-id 27
$ws.Cells.Item(215, 1).Value = 46007
$ws.Cells.Item(215, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(215, 3).Value = 300.29500000000002
$ws.Cells.Item(215, 4).Value = 305.34400000000005
$ws.Cells.Item(215, 5).Value = 181.21899999999999
$ws.Cells.Item(215, 6).Value = 51.633000000000003
$ws.Cells.Item(215, 7).Value = 136.82
$ws.Cells.Item(215, 8).Value = 155.149
$ws.Cells.Item(215, 9).Value = 57.414999999999999
$ws.Cells.Item(215, 10).Value = 118.492
$ws.Cells.Item(215, 11).Value = 300.14400000000001
$ws.Cells.Item(215, 12).Value = 214.71899999999999
$ws.Cells.Item(215, 13).Value = 163.04399999999998
$ws.Cells.Item(215, 14).Value = 242.07400000000001
$ws.Cells.Item(215, 15).Value = 359.23199999999997
$ws.Cells.Item(215, 16).Value = 531.774
$ws.Cells.Item(215, 17).Value = 322.73500000000001
$ws.Cells.Item(215, 18).Value = 382.74999999999994
$ws.Cells.Item(215, 19).Value = 210.345
$ws.Cells.Item(215, 20).Value = 131.49600000000001
$ws.Cells.Item(215, 21).Value = 5.7969999999999997
$ws.Cells.Item(215, 22).Value = 117.97500000000002
$ws.Cells.Item(215, 23).Value = 41.796999999999997
$ws.Cells.Item(215, 24).Value = 0
$ws.Cells.Item(215, 25).Value = 9.9990000000000006
$ws.Cells.Item(215, 26).Value = 26.243000000000002

# Row 216: date serial 46008, station string-id 26
$ws.Cells.Item(216, 1).Value = 46008
$ws.Cells.Item(216, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(216, 3).Value = 628.17699999999991
$ws.Cells.Item(216, 4).Value = 836.79600000000005
$ws.Cells.Item(216, 5).Value = 432.36399999999998
$ws.Cells.Item(216, 6).Value = 251.66
$ws.Cells.Item(216, 7).Value = 350.26900000000001
$ws.Cells.Item(216, 8).Value = 627.99300000000017
$ws.Cells.Item(216, 9).Value = 350.75900000000007
$ws.Cells.Item(216, 10).Value = 92.686000000000007
$ws.Cells.Item(216, 11).Value = 206.798
$ws.Cells.Item(216, 12).Value = 131.328
$ws.Cells.Item(216, 13).Value = 154.80800000000002
$ws.Cells.Item(216, 14).Value = 249.24600000000001
$ws.Cells.Item(216, 15).Value = 572.50900000000013
$ws.Cells.Item(216, 16).Value = 1241.4489999999998
$ws.Cells.Item(216, 17).Value = 362.82
$ws.Cells.Item(216, 18).Value = 309.15699999999998
$ws.Cells.Item(216, 19).Value = 413.01299999999998
$ws.Cells.Item(216, 20).Value = 100.22999999999999
$ws.Cells.Item(216, 21).Value = 117.482
$ws.Cells.Item(216, 22).Value = 14.8
$ws.Cells.Item(216, 23).Value = 105.06
$ws.Cells.Item(216, 24).Value = 62.260000000000005
$ws.Cells.Item(216, 25).Value = 77.910000000000011
$ws.Cells.Item(216, 26).Value = 113.00699999999999

# Row 217: date serial 46008, station string-id 27
$ws.Cells.Item(217, 1).Value = 46008
$ws.Cells.Item(217, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(217, 3).Value = 451.41900000000004
$ws.Cells.Item(217, 4).Value = 503.24100000000004
$ws.Cells.Item(217, 5).Value = 119.28100000000001
$ws.Cells.Item(217, 6).Value = 20.484000000000002
$ws.Cells.Item(217, 7).Value = 15.797000000000001
$ws.Cells.Item(217, 8).Value = 99.397999999999996
$ws.Cells.Item(217, 9).Value = 330.43699999999995
$ws.Cells.Item(217, 10).Value = 159.108
$ws.Cells.Item(217, 11).Value = 130.17099999999999
$ws.Cells.Item(217, 12).Value = 321.565
$ws.Cells.Item(217, 13).Value = 167.048
$ws.Cells.Item(217, 14).Value = 186.29499999999996
$ws.Cells.Item(217, 15).Value = 294.298
$ws.Cells.Item(217, 16).Value = 597.85599999999988
$ws.Cells.Item(217, 17).Value = 377.92500000000001
$ws.Cells.Item(217, 18).Value = 182.834
$ws.Cells.Item(217, 19).Value = 160.904
$ws.Cells.Item(217, 20).Value = 66.114999999999995
$ws.Cells.Item(217, 21).Value = 66.164999999999992
$ws.Cells.Item(217, 22).Value = 55.460999999999999
$ws.Cells.Item(217, 23).Value = 0
$ws.Cells.Item(217, 24).Value = 94.836999999999989
$ws.Cells.Item(217, 25).Value = 18.171999999999997
$ws.Cells.Item(217, 26).Value = 19.501000000000001

# Row 218: date serial 46009, station string-id 26
$ws.Cells.Item(218, 1).Value = 46009
$ws.Cells.Item(218, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(218, 3).Value = 565.22400000000005
$ws.Cells.Item(218, 4).Value = 891.78099999999995
$ws.Cells.Item(218, 5).Value = 254.97
$ws.Cells.Item(218, 6).Value = 432.74
$ws.Cells.Item(218, 7).Value = 325.45400000000001
$ws.Cells.Item(218, 8).Value = 520.52099999999996
$ws.Cells.Item(218, 9).Value = 419.73
$ws.Cells.Item(218, 10).Value = 196.79899999999998
$ws.Cells.Item(218, 11).Value = 25.880000000000003
$ws.Cells.Item(218, 12).Value = 158.80000000000001
$ws.Cells.Item(218, 13).Value = 227.78299999999999
$ws.Cells.Item(218, 14).Value = 187.27799999999999
$ws.Cells.Item(218, 15).Value = 700.53400000000011
$ws.Cells.Item(218, 16).Value = 1594.0129999999999
$ws.Cells.Item(218, 17).Value = 413.12
$ws.Cells.Item(218, 18).Value = 476.01300000000003
$ws.Cells.Item(218, 19).Value = 368.11800000000005
$ws.Cells.Item(218, 20).Value = 154.535
$ws.Cells.Item(218, 21).Value = 20.57
$ws.Cells.Item(218, 22).Value = 27.531999999999996
$ws.Cells.Item(218, 23).Value = 230.04999999999998
$ws.Cells.Item(218, 24).Value = 66.62
$ws.Cells.Item(218, 25).Value = 58.73
$ws.Cells.Item(218, 26).Value = 45.8

# Row 219: date serial 46009, station string-id 27
$ws.Cells.Item(219, 1).Value = 46009
$ws.Cells.Item(219, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(219, 3).Value = 312.73700000000002
$ws.Cells.Item(219, 4).Value = 439.03300000000002
$ws.Cells.Item(219, 5).Value = 53.146000000000001
$ws.Cells.Item(219, 6).Value = 131.03100000000001
$ws.Cells.Item(219, 7).Value = 20.609000000000002
$ws.Cells.Item(219, 8).Value = 91.194999999999993
$ws.Cells.Item(219, 9).Value = 263.75299999999999
$ws.Cells.Item(219, 10).Value = 177.23400000000001
$ws.Cells.Item(219, 11).Value = 98.015000000000001
$ws.Cells.Item(219, 12).Value = 199.70200000000003
$ws.Cells.Item(219, 13).Value = 242.45800000000003
$ws.Cells.Item(219, 14).Value = 297.16899999999993
$ws.Cells.Item(219, 15).Value = 320.98999999999995
$ws.Cells.Item(219, 16).Value = 620.17100000000005
$ws.Cells.Item(219, 17).Value = 249.60900000000004
$ws.Cells.Item(219, 18).Value = 95.319000000000017
$ws.Cells.Item(219, 19).Value = 256.31399999999996
$ws.Cells.Item(219, 20).Value = 328.18799999999999
$ws.Cells.Item(219, 21).Value = 51.202999999999996
$ws.Cells.Item(219, 22).Value = 0
$ws.Cells.Item(219, 23).Value = 0
$ws.Cells.Item(219, 24).Value = 67.262
$ws.Cells.Item(219, 25).Value = 0
$ws.Cells.Item(219, 26).Value = 23.422000000000001

# Restore the cursor/selection position recorded in the commit
$ws.Range("D226").Select()
